$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback record being appended (row 4) across all three sheets.
# ---------------------------------------------------------------------------
$mdName      = "d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md"
$mdPath      = "e2e\d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md"
$statusSync  = "Handed back: in sync with en-US"
$hoDate      = "2016-08-28 14:44:17"

$linkColor = 0xED9564   # BGR for RGB(0x64,0x95,0xED) -> matches existing HyperLink font color
$dateFmt   = "yyyy-mm-dd HH:mm:ss"

function Format-AsHyperlinkText($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $linkColor
}

function Format-AsDate($rng) {
    $rng.NumberFormat = $dateFmt
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $mdPath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusSync
$wsOverview.Range("F4").Value = $statusSync
$wsOverview.Range("G4").Value = $hoDate

Format-AsHyperlinkText $wsOverview.Range("B4")
Format-AsDate $wsOverview.Range("G4")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7590e1eae5ae6add05f747bc2086873722ae62/e2e/d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md", "", "", $mdPath) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = $mdName
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.a7590e1eae5ae6add05f747bc2086873722ae62e.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-08-28 14:44:13"
$wsZhCn.Range("I4").Value = $mdName
$wsZhCn.Range("J4").Value = "d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.a7590e1eae5ae6add05f747bc2086873722ae62e.zh-cn.xlf"
$wsZhCn.Range("K4").Value = "2016-08-28 14:44:30"
$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

Format-AsHyperlinkText $wsZhCn.Range("A4")
Format-AsHyperlinkText $wsZhCn.Range("I4")
Format-AsDate $wsZhCn.Range("H4")
Format-AsDate $wsZhCn.Range("K4")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7590e1eae5ae6add05f747bc2086873722ae62/e2e/d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md", "", "", $mdName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/a7590e1eae5ae6add05f747bc2086873722ae62/e2e/d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md", "", "", $mdName) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = $mdName
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.a7590e1eae5ae6add05f747bc2086873722ae62e.de-de.xlf"
$wsDeDe.Range("H4").Value = $hoDate
$wsDeDe.Range("I4").Value = $mdName
$wsDeDe.Range("J4").Value = "d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.a7590e1eae5ae6add05f747bc2086873722ae62e.de-de.xlf"
$wsDeDe.Range("K4").Value = "2016-08-28 14:44:37"
$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

Format-AsHyperlinkText $wsDeDe.Range("A4")
Format-AsHyperlinkText $wsDeDe.Range("I4")
Format-AsDate $wsDeDe.Range("H4")
Format-AsDate $wsDeDe.Range("K4")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a7590e1eae5ae6add05f747bc2086873722ae62/e2e/d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md", "", "", $mdName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a7590e1eae5ae6add05f747bc2086873722ae62/e2e/d3b7e7d8-23aa-4034-968e-a7a9b2e87c8e.md", "", "", $mdName) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Output "Handback report row appended to Overview, zh-cn and de-de sheets."
